$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting existing A:E -> B:F
$ws.Columns("A:A").Insert()

# New column A header
$ws.Range("A1").Value = "ID"

# New column A row labels (IDs) for rows 2-25
$ids = @(
  "Hb 2", "Hb 3", "S 24", "S 28", "Hb 107", "Hb 66", "Hb 69", "Hb 95",
  "Hb 99", "Hb 92", "Hb 40", "Hb 41", "S 11", "Hb 57", "S 21", "S 22",
  "S 3", "S 4", "S 5", "Hb 74", "Hb 79", "Hb 32", "S 15", "S 16"
)

for ($i = 0; $i -lt $ids.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 1).Value = $ids[$i]
}

# Match the header style (A1) to the rest of row 1 (bordered/bold/centered header)
$ws.Range("B1:F1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
